# New crime data collected — weekly CompStat update for the 115th Precinct.
#
# Updates the "Volume/Number" and "Report Covering the Week" rich-text
# captions (in-place character replacement so existing run formatting is
# preserved), plus the Week-to-Date / 28-Day / Year-to-Date / 2-Year crime
# figures in rows 14-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header captions: edit only the affected run of text, leaving the rest
# of the rich-text cell (and its per-run formatting) untouched.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  47" -> "...48"
$a8 = $ws.Range("A8")
$a8text = $a8.Value()
$pos = $a8text.IndexOf("47") + 1
$a8.Characters($pos, 2).Text = "48"

# C9: "Report Covering the Week  11/20/2023  Through  11/26/2023"
#  -> "...11/27/2023  Through  12/3/2023"
$c9 = $ws.Range("C9")
$c9text = $c9.Value()
$pos1 = $c9text.IndexOf("11/20/2023") + 1
$c9.Characters($pos1, 10).Text = "11/27/2023"
$c9text2 = $c9.Value()
$pos2 = $c9text2.IndexOf("11/26/2023") + 1
$c9.Characters($pos2, 10).Text = "12/3/2023"

# ---------------------------------------------------------------------
# Crime-complaints table (rows 14-30): updated counts / percentages.
# ---------------------------------------------------------------------

# Row 14
$ws.Range("M14").Value = 0

# Row 15
$ws.Range("M15").Value = -20.689655172413
$ws.Range("N15").Value = -30.303030303030

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 328
$ws.Range("J16").Value = 290
$ws.Range("K16").Value = 13.103448275862
$ws.Range("L16").Value = 43.859649122807
$ws.Range("M16").Value = 1.547987616099
$ws.Range("N16").Value = -72.959604286892

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 30
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 441
$ws.Range("J17").Value = 384
$ws.Range("K17").Value = 14.84375
$ws.Range("L17").Value = 22.5
$ws.Range("M17").Value = 45.065789473684
$ws.Range("N17").Value = 10.25

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 130
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = -7.142857142857
$ws.Range("L18").Value = -1.515151515151
$ws.Range("M18").Value = -52.205882352941
$ws.Range("N18").Value = -92.923244420250

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 790
$ws.Range("J19").Value = 909
$ws.Range("K19").Value = -13.091309130913
$ws.Range("L19").Value = 21.351766513056
$ws.Range("M19").Value = 69.527896995708
$ws.Range("N19").Value = -41.263940520446

# Row 20
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 36
$ws.Range("I20").Value = 324
$ws.Range("J20").Value = 292
$ws.Range("K20").Value = 10.958904109589
$ws.Range("L20").Value = 69.633507853403
$ws.Range("M20").Value = 45.945945945945
$ws.Range("N20").Value = -84.527220630372

# Row 21
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 165
$ws.Range("G21").Value = 162
$ws.Range("H21").Value = 1.851851851851
$ws.Range("I21").Value = 2039
$ws.Range("J21").Value = 2054
$ws.Range("K21").Value = -0.730282375851
$ws.Range("L21").Value = 27.676894176581
$ws.Range("M21").Value = 25.941939468807
$ws.Range("N21").Value = -70.615362444156

# Row 22 (G22/H22 flip from numeric 0-pct cells to the "N/A" text markers
# already used elsewhere in the table, so copy format+value from a cell
# that already carries that exact text/style combination).
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 7
$ws.Range("D14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("I22").Value = 83
$ws.Range("K22").Value = 43.103448275862
$ws.Range("L22").Value = 196.428571428571
$ws.Range("M22").Value = 196.428571428571

# Row 24
$ws.Range("C24").Value = 54
$ws.Range("E24").Value = 35
$ws.Range("F24").Value = 169
$ws.Range("G24").Value = 191
$ws.Range("H24").Value = -11.518324607329
$ws.Range("I24").Value = 1853
$ws.Range("J24").Value = 1810
$ws.Range("K24").Value = 2.375690607734
$ws.Range("L24").Value = 45.790715971675
$ws.Range("M24").Value = 76.644423260247

# Row 25
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -23.809523809523
$ws.Range("F25").Value = 80
$ws.Range("G25").Value = 75
$ws.Range("H25").Value = 6.666666666666
$ws.Range("I25").Value = 911
$ws.Range("J25").Value = 843
$ws.Range("K25").Value = 8.066429418742
$ws.Range("L25").Value = 15.462610899873
$ws.Range("M25").Value = 5.561993047508

# Row 26
$ws.Range("L26").Value = -15.094339622641

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -36.363636363636
$ws.Range("I27").Value = 135
$ws.Range("J27").Value = 105
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 37.755102040816

# Row 28
$ws.Range("N28").Value = -90.909090909090

# Row 29
$ws.Range("N29").Value = -92.156862745098

# Row 30 (C30 flips from numeric 1 to the "N/A" text marker)
$ws.Range("D14").Copy($ws.Range("C30"))
